# Add Q3-2022 reporting period: a new "2022-Q3" detail sheet (placed right
# after "总计" / before "2022-Q2") plus a new summary row on "总计" for
# that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" sheet before the existing "2022-Q2" sheet.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Bold / bordered / centered look matching the other quarter sheets'
# header row and index column.
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$rows = @(
    @("550015", "中信保诚至远动力混合A", "16.29", "75.56", "3.18", "0.5180", 6),
    @("009913", "中信保诚成长动力混合A", "5.03", "74.01", "3.65", "0.1836", 5),
    @("550016", "中信保诚至远动力混合C", "2.30", "75.56", "3.18", "0.0731", 6),
    @("011351", "金鹰年年邮益一年持有期混合A", "3.43", "34.33", "0.74", "0.0254", 8),
    @("014831", "兴银中证1000指数增强A", "1.37", "83.33", "1.09", "0.0149", 3),
    @("014832", "兴银中证1000指数增强C", "0.90", "83.33", "1.09", "0.0098", 3),
    @("014282", "中信保诚成长动力混合C", "0.07", "74.01", "3.65", "0.0026", 5),
    @("011352", "金鹰年年邮益一年持有期混合C", "0.27", "34.33", "0.74", "0.0020", 8),
    @("014677", "中信保诚至远动力混合E", "0.01", "75.56", "3.18", "0.0003", 6)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = ($r - 2)
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

$indexRange = $newSheet.Range("A2:A10")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2) Insert a new row 2 on "总计" for the 2022-Q3 summary, shifting the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Row 6 is brand new on this sheet - prime its column-A cell with the
# same style as the row above it before filling in the shifted value.
$totalSheet.Cells.Item(5, 1).Copy($totalSheet.Cells.Item(6, 1))

# Shift quarter rows 2-5 down into rows 3-6 (bottom-up so sources aren't
# overwritten before they're read).
for ($row = 5; $row -ge 2; $row--) {
    $totalSheet.Cells.Item($row + 1, 2).Value = $totalSheet.Cells.Item($row, 2).Value2
    $totalSheet.Cells.Item($row + 1, 3).Value = $totalSheet.Cells.Item($row, 3).Value2
    $totalSheet.Cells.Item($row + 1, 4).Value = $totalSheet.Cells.Item($row, 4).Value2
}

# Renumber the index column (A), 0-based.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4

# Write the new 2022-Q3 summary row.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 9
$totalSheet.Cells.Item(2, 4).Value = 0.83
